# "Generate Report for Handoff"
# The localization run that produced this report has moved on from
# "In Translation" to "Ready for handoff" -- refresh the status text and
# the "latest xliff generated" timestamps on all three sheets, then widen
# the now-longer status columns to fit the new text (mirrors what the
# report generator's own AutoFit pass does).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- status text: "In Translation" -> "Ready for handoff" ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value     = "Ready for handoff"
$dede.Range("C2").Value     = "Ready for handoff"

# --- refreshed "latest xliff generate / handoff" timestamps ---
$overview.Range("G2").Value = "2016-09-03 02:42:39"
$zhcn.Range("H2").Value     = "2016-09-03 02:42:35"
$dede.Range("H2").Value     = "2016-09-03 02:42:39"

# --- widen the status columns to fit "Ready for handoff" ---
$overview.Range("E1:F1").ColumnWidth = 16.333333333333336
$zhcn.Range("C1").ColumnWidth        = 16.333333333333336
$dede.Range("C1").ColumnWidth        = 16.333333333333336
